$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.14"
$ws.Range("D3").Value = "'23.30"
$ws.Range("D4").Value = "'5.417"
$ws.Range("D5").Value = "'0.05982"
$ws.Range("D7").Value = "'6.529"
$ws.Range("D8").Value = "'0.8147"
$ws.Range("D9").Value = "'0.9199"
$ws.Range("D10").Value = "'0.1409"
$ws.Range("D11").Value = "'0.07406"
$ws.Range("D12").Value = "'0.03236"
$ws.Range("D14").Value = "'0.09342"
$ws.Range("D15").Value = "'3.847"
$ws.Range("D16").Value = "'0.001571"
$ws.Range("D17").Value = "'0.04677"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005939"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "'0.006087"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Value = "'0.005011"
$ws.Range("E20").Value = "19HotbitTokenHTBBestin24h"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").Value = "'0.0009804"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "'0.00007799"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'3.622"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "'2.130"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("D40").Value = "'0.03923"
$ws.Range("D41").Value = "'0.006218"
$ws.Range("D43").Value = "'0.002620"
$ws.Range("D44").Value = "'0.007083"
$ws.Range("D45").Value = "'0.00005244"
$ws.Range("D48").Value = "'0.9099"

Write-Host "Applied 48 cell updates"
